# Apply cryptocurrency price/volume updates per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '68.028.75'
$ws.Range('E2').Value = '  +0.57%  '

# Row 3
$ws.Range('D3').Value = '3.780.63'
$ws.Range('E3').Value = '  -0.39%  '

# Row 4
$ws.Range('E4').Value = '  -0.13%  '

# Row 5
$ws.Range('D5').Value = "'599.29"
$ws.Range('E5').Value = '  +0.46%  '

# Row 6
$ws.Range('D6').Value = "'163.26"
$ws.Range('E6').Value = '  -2.32%  '

# Row 7
$ws.Range('E7').Value = '  -0.03%  '

# Row 8
$ws.Range('E8').Value = '  -0.98%  '

# Row 9
$ws.Range('E9').Value = '  -1.64%  '

# Row 10
$ws.Range('D10').Value = "'0.447"
$ws.Range('E10').Value = '  -0.42%  '

# Row 11
$ws.Range('D11').Value = "'6.59"
$ws.Range('E11').Value = '  +4.43%  '

# Row 12
$ws.Range('E12').Value = '  -2.63%  '

# Row 13
$ws.Range('D13').Value = "'35.37"
$ws.Range('E13').Value = '  -1.29%  '

# Row 14
$ws.Range('D14').Value = '4.414.39'
$ws.Range('E14').Value = '  -0.58%  '

# Row 15
$ws.Range('D15').Value = '3.808.15'
$ws.Range('E15').Value = '  -0.72%  '

# Row 16
$ws.Range('D16').Value = '67.953.85'
$ws.Range('E16').Value = '  +0.40%  '

# Row 17
$ws.Range('D17').Value = "'18.25"
$ws.Range('E17').Value = '  -1.65%  '

# Row 18
$ws.Range('E18').Value = '  +2.00%  '

# Row 19
$ws.Range('D19').Value = "'7.01"
$ws.Range('E19').Value = '  -0.85%  '

# Row 20
$ws.Range('D20').Value = "'457.82"
$ws.Range('E20').Value = '  -0.72%  '

# Row 21
$ws.Range('D21').Value = "'9.59"
$ws.Range('E21').Value = '  -3.36%  '

# Row 22
$ws.Range('D22').Value = "'0.695"
$ws.Range('E22').Value = '  -0.66%  '

# Row 23
$ws.Range('D23').Value = "'82.75"
$ws.Range('E23').Value = '  -0.80%  '

# Row 24
$ws.Range('E24').Value = '  -5.92%  '

# Row 25
$ws.Range('D25').Value = "'11.93"
$ws.Range('E25').Value = '  -1.29%  '

# Row 26
$ws.Range('D26').Value = "'2.08"
$ws.Range('E26').Value = '  -0.80%  '

# Row 27
$ws.Range('E27').Value = '  -0.02%  '

# Row 28
$ws.Range('D28').Value = "'9.91"
$ws.Range('E28').Value = '  -0.93%  '

# Row 29
$ws.Range('D29').Value = '3.928.98'
$ws.Range('E29').Value = '  -0.32%  '

# Row 30
$ws.Range('E30').Value = '  -0.68%  '

# Row 31
$ws.Range('E31').Value = '  -1.00%  '

# Row 32
$ws.Range('D32').Value = "'2.56"
$ws.Range('E32').Value = '  -7.56%  '

# Row 33
$ws.Range('D33').Value = "'28.99"
$ws.Range('E33').Value = '  -1.83%  '

# Row 34
$ws.Range('E34').Value = '  -0.01%  '

# Row 35
$ws.Range('D35').Value = "'8.93"
$ws.Range('E35').Value = '  -1.36%  '

# Row 36
$ws.Range('D36').Value = "'0.0990"
$ws.Range('E36').Value = '  -0.97%  '

# Row 37
$ws.Range('D37').Value = "'0.142"
$ws.Range('E37').Value = '  +2.98%  '

# Row 38
$ws.Range('E38').Value = '  +0.35%  '

# Row 39
$ws.Range('D39').Value = "'0.980"
$ws.Range('E39').Value = '  -1.62%  '

# Row 40
$ws.Range('E40').Value = '  -5.89%  '

# Row 41
$ws.Range('E41').Value = '  -0.08%  '

# Row 43
$ws.Range('B43').Value = 'OKB'
$ws.Range('C43').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D43').Value = "'47.28"
$ws.Range('E43').Value = '  -1.66%  '

# Row 44
$ws.Range('B44').Value = 'Arweave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D44').Value = "'43.31"
$ws.Range('E44').Value = '  +0.84%  '

# Row 45
$ws.Range('D45').Value = "'152.68"
$ws.Range('E45').Value = '  +3.14%  '

# Row 46
$ws.Range('E46').Value = '  -1.68%  '

# Row 47
$ws.Range('B47').Value = 'Cosmos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D47').Value = "'8.29"
$ws.Range('E47').Value = '  -0.29%  '

# Row 48
$ws.Range('B48').Value = 'ONDO'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D48').Value = "'1.36"
$ws.Range('E48').Value = '  +1.19%  '

# Row 49
$ws.Range('D49').Value = "'1.85"
$ws.Range('E49').Value = '  +0.41%  '

# Row 50
$ws.Range('D50').Value = "'387.18"
$ws.Range('E50').Value = '  -2.20%  '

# Row 51
$ws.Range('D51').Value = "'26.32"
$ws.Range('E51').Value = '  -3.69%  '
